# Filtering of sources added
# Replace "Backend Apis" / "Api Controllers need to be moved to a separate project"
# with "Notification" / "Toastr Integration" at row 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B25").Value = "Notification"
$ws.Range("C25").Value = "Toastr Integration"

# Update the active selection to reflect the edited row (cosmetic view state change)
$ws.Range("C26").Select()
